# deliveries.xlsx: remove unused fields (authorization_name, authorization_description)
# and add distribution_name to the "model" (entitlements/deliveries) sheet; also reorder
# item_pack_name/item_pack_description. Also widen settings!A and move some selections.

$wb = $excel.ActiveWorkbook

$modelSheet = $wb.Worksheets.Item("model")

# Row 4 currently holds "authorization_name" -> rename field to "distribution_name"
$modelSheet.Range("B4").Value = "distribution_name"

# Row 5 currently holds "authorization_description" -> remove this row entirely,
# shifting everything below it up by one.
$modelSheet.Rows(5).Delete()

# After the delete, the old rows:
#   11 item_pack_description
#   12 item_pack_name
# are now:
#   10 item_pack_description
#   11 item_pack_name
# Swap them so item_pack_name comes before item_pack_description.
$tmp = $modelSheet.Range("B10").Value()
$modelSheet.Range("B10").Value = $modelSheet.Range("B11").Value()
$modelSheet.Range("B11").Value = $tmp

# "settings" sheet: widen column A and move the selection.
$settingsSheet = $wb.Worksheets.Item("settings")
$settingsSheet.Columns("A").ColumnWidth = 21.7
$settingsSheet.Range("C14").Select()

# Restore "model" as the active sheet/tab and move its selection, matching the
# original workbook (model stays the active tab).
$modelSheet.Activate()
$modelSheet.Range("B7").Select()
